# Apply the "Metadata" + "Elements" sheet updates described by the diff.
$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> Alvearie Team
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" row -> becomes "Jurisdiction" / "United States of America"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was the second (now redundant) "Contact" row -> delete it entirely,
# shifting "Description" and everything below up by one row.
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elem = $wb.Worksheets.Item("Elements")

# Root element row: Short / Definition get the real extension text
$elem.Cells.Item(2, 11).Value = "Communication Care Gap"
$elem.Cells.Item(2, 12).Value = "ID of care gap for which we are communicating"
